$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 168471
$ws.Range("C4").Value = 159334
$ws.Range("C5").Value = 9137
$ws.Range("C8").Value = 65.56999999999999
